$d = $word.ActiveDocument

# --- 1. Target the first paragraph (the **ID__...__ID** placeholder paragraph) ---
$p1 = $d.Paragraphs.Item(1)
$pf = $p1.Range.ParagraphFormat

# Add a paragraph border (top/left/bottom/right) whose only attribute is the
# "space" (distance from text) = 5, without turning on a visible line style.
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$pf.LeftIndent = 11.25

# --- 2. Update the placeholder ID text and drop the trailing space run ---
$oldId = "**ID__AFFARS_mp_5315_3_topic_4__ID**"
$newId = "**ID__AFFARS_MP_5315_3_1_2__ID**"

$found = $d.Content.Find.Execute($oldId, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $idRange = $d.Content
    $idRange.Text = $newId

    $afterStart = $idRange.Start + $newId.Length
    $afterRange = $d.Range($afterStart, $afterStart + 1)
    if ($afterRange.Text -eq " ") {
        $afterRange.Delete()
    }
}
